$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "Item Count" column (M) ---
# Pick up the same formatting already used in column L, then fill in the
# header text and a count of 1 for every existing food item row.
$ws.Range("L1:L11").Copy()
$ws.Range("M1:M11").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("M1").Value = "Item Count"
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 13).Value = 1
}

# Row heights for the data rows recalculate (rows 1-2 keep their explicit
# heights, the rest revert to auto height) once the new column is in place.
$ws.Rows("3:11").AutoFit()

# --- Extend the table's formatting further down the sheet ---
$ws.Range("A2:L11").Copy()
$ws.Range("A12:L25").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A12:L25").ClearContents()

# Row 13 ends up blank/unformatted.
$ws.Range("A13:L13").Clear()

# Row 14's first cell picks up the header style instead of the body style.
$ws.Range("B1").Copy()
$ws.Range("A14").PasteSpecial(-4122)  # xlPasteFormats

$ws.Application.CutCopyMode = $false

# Final selection left by the editing session.
$ws.Range("G20").Select()
